$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.898.30"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "'1.635.13"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'216.36"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "'0.5072"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.2578"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "'0.06362"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'4.260"
$ws.Range("D13").Value = "'1.635.01"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "'0.5517"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'0.0₅7710"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").Value = "'64.03"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "'25.918.14"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'4.448"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'194.43"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'9.896"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'1.908"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'142.66"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").Value = "'0.1236"
$ws.Range("E26").Value = "  +5.88%  "
$ws.Range("D27").Value = "'6.823"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "'0.04866"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D31").Value = "'3.251"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'3.189"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "'1.544"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "'2.374"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "'0.9053"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'2.569"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'0.5494"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "'1.122.89"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D40").Value = "'1.002"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "'5.578"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'0.8055"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").Value = "'97.50"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").Value = "'0.0₈121"
$ws.Range("E44").Value = "  -6.43%  "
$ws.Range("D45").Value = "'1.774.50"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'0.4459"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("D47").Value = "'54.89"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'0.9963"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").Value = "'0.05152"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").Value = "'7.542"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("D51").Value = "'1.005"
$ws.Range("E51").Value = "  +0.02%  "
